$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set C5 to DONE, matching the other "DONE" entries - this will also
# update the dependent shared formula in D5 (and the SUM total in D13).
$ws.Range("C5").Value = "DONE"

# Update the selected cell/range as recorded in the sheet view.
$ws.Range("F8").Select()
